$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = 25.46092070582212
$ws.Range("D2").Value = 48.89928110391109
$ws.Range("E2").Value = 29.36823407648259
$ws.Range("F2").Value = 100
$ws.Range("G2").Value = 79.76887265481093
$ws.Range("C3").Value = 35.47070001895196
$ws.Range("D3").Value = 13.43946833874935
$ws.Range("E3").Value = 26.33800329683886
$ws.Range("F3").Value = 44.82603157826325
$ws.Range("G3").Value = 76.34167738253856
$ws.Range("C4").Value = 62.65030358103355
$ws.Range("D4").Value = 20.328262771856
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 45.0693541582294
$ws.Range("G4").Value = 100
$ws.Range("C5").Value = 53.16831676755705
$ws.Range("D5").Value = 15.38001080534552
$ws.Range("E5").Value = 62.04430198085566
$ws.Range("F5").Value = 43.87679514285356
$ws.Range("G5").Value = 57.66466763358127
$ws.Range("C6").Value = 56.72933313736404
$ws.Range("D6").Value = 80.46100733464016
$ws.Range("E6").Value = 66.29812480115532
$ws.Range("F6").Value = 100
$ws.Range("G6").Value = 29.89297777248765
$ws.Range("C7").Value = 79.56627411049811
$ws.Range("D7").Value = 19.48720233820424
$ws.Range("E7").Value = 94.7708932802425
$ws.Range("F7").Value = 41.70673899721343
$ws.Range("G7").Value = 76.58343204936317
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = 14.68810980988851
$ws.Range("E8").Value = 75.71909102655357
$ws.Range("F8").Value = 79.27647896967049
$ws.Range("G8").Value = 100
$ws.Range("C9").Value = 61.77368338896039
$ws.Range("D9").Value = 72.70511039761431
$ws.Range("E9").Value = 100
$ws.Range("F9").Value = 76.36980188576024
$ws.Range("G9").Value = 43.46939083426755
$ws.Range("C10").Value = 62.29318466346002
$ws.Range("D10").Value = 51.51598588623294
$ws.Range("E10").Value = 90.50399979915296
$ws.Range("F10").Value = 45.94876242539156
$ws.Range("G10").Value = 100
$ws.Range("C11").Value = 48.36269504309307
$ws.Range("D11").Value = 24.49163739036445
$ws.Range("E11").Value = 39.69592121582553
$ws.Range("F11").Value = 100
$ws.Range("G11").Value = 73.72837648363699
$ws.Range("C12").Value = 26.21767906357819
$ws.Range("D12").Value = 43.97070949561198
$ws.Range("E12").Value = 88.21132617235055
$ws.Range("F12").Value = 83.84160136738717
$ws.Range("G12").Value = 78.14059510703973
$ws.Range("C13").Value = 100
$ws.Range("D13").Value = 76.01442758840712
$ws.Range("E13").Value = 76.24009844570018
$ws.Range("F13").Value = 37.17206092743768
$ws.Range("G13").Value = 100
$ws.Range("C14").Value = 90.68989183504158
$ws.Range("D14").Value = 58.30508627637893
$ws.Range("E14").Value = 100
$ws.Range("F14").Value = 86.19931411257902
$ws.Range("G14").Value = 64.38484309842855
$ws.Range("C15").Value = 30.81115020868923
$ws.Range("D15").Value = 48.02463459268451
$ws.Range("E15").Value = 52.44895385076189
$ws.Range("F15").Value = 100
$ws.Range("G15").Value = 95.65773559565649
$ws.Range("C16").Value = 55.44047744661044
$ws.Range("D16").Value = 41.027841706961
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 92.67614514625859
$ws.Range("G16").Value = 100
$ws.Range("C17").Value = 56.56134163456932
$ws.Range("D17").Value = 30.29428735738359
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 15.46691612318783
$ws.Range("G17").Value = 92.36382799119454
$ws.Range("C18").Value = 65.53142484668035
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 99.80490147995826
$ws.Range("G18").Value = 81.5009264733345
$ws.Range("C19").Value = 38.10080483931941
$ws.Range("D19").Value = 8.560731798254054
$ws.Range("E19").Value = 99.20256804176583
$ws.Range("F19").Value = 55.34113756524128
$ws.Range("G19").Value = 40.00522405864503
$ws.Range("C20").Value = 28.00689350224013
$ws.Range("D20").Value = 62.2724052047541
$ws.Range("E20").Value = 78.96361241406994
$ws.Range("F20").Value = 100
$ws.Range("G20").Value = 96.0828309977109
$ws.Range("C21").Value = 95.04943788441064
$ws.Range("D21").Value = 68.77969285782628
$ws.Range("E21").Value = 30.38223535434911
$ws.Range("F21").Value = 48.66982999772584
$ws.Range("G21").Value = 38.51221430445766
